# Update contestant seating assignments with new identifiers
# Swap the ID (column A) and ContestantID (column C) values between
# row 2 and row 3 of the "Seat Assignments" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Seat Assignments")

$ws.Range("A2").Value = "91e25164-6f67-42f7-b978-9132a406c060"
$ws.Range("C2").Value = "28603f95-d5f6-47ab-88c4-0d79742a6b02"

$ws.Range("A3").Value = "ff87f03b-8891-4bb6-ac5c-a510d216fdd6"
$ws.Range("C3").Value = "d698b1de-6641-45c6-aa63-f577d2b634bb"
